$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data block (row 578), pushing the
# existing rows (old 578-606) down to 582-610.
$ws.Rows("578:581").Insert()

# Shared/constant values for this market block (identical across the sheet).
$mercadoId = 2
$mercado = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102005
$categoria = "Naranja"
$unidad = "`$/bins (400 kilos)"
$origen = "Provincia de Limarí"
$kgUnidad = 400

# New weekly rows (fecha 45267) appended for this product.
$newRows = @(
    @{ Row = 578; Fecha = 45267; Variedad = "Lane Late"; Calidad = "Primera"; Volumen = 20; PMin = 190000; PMax = 200000; PProm = 195000; PKg = 488 },
    @{ Row = 579; Fecha = 45267; Variedad = "Lane Late"; Calidad = "Segunda"; Volumen = 16; PMin = 150000; PMax = 160000; PProm = 155000; PKg = 388 },
    @{ Row = 580; Fecha = 45267; Variedad = "Valencia";  Calidad = "Primera"; Volumen = 20; PMin = 200000; PMax = 210000; PProm = 205000; PKg = 512 },
    @{ Row = 581; Fecha = 45267; Variedad = "Valencia";  Calidad = "Segunda"; Volumen = 16; PMin = 160000; PMax = 170000; PProm = 165000; PKg = 412 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
